# Apply weekly crime-data refresh for the 34th Precinct CompStat report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and reporting week ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Crime complaint statistics table (rows 14-29) ---
# Row 14
$c = $ws.Range("M14")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -100.0

# Row 16
$c = $ws.Range("C16")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("D16")
$c.NumberFormat = "#,##0"
$c.Value = 3.0
$c = $ws.Range("E16")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -66.666666666666
$c = $ws.Range("G16")
$c.NumberFormat = "#,##0"
$c.Value = 13.0
$c = $ws.Range("H16")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 15.384615384615
$c = $ws.Range("I16")
$c.NumberFormat = "#,##0"
$c.Value = 81.0
$c = $ws.Range("J16")
$c.NumberFormat = "#,##0"
$c.Value = 77.0
$c = $ws.Range("K16")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 5.194805194805
$c = $ws.Range("L16")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 50.0
$c = $ws.Range("M16")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -22.115384615384
$c = $ws.Range("N16")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -84.65909090909

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "#,##0"
$c.Value = 8.0
$c = $ws.Range("E17")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -50.0
$c = $ws.Range("G17")
$c.NumberFormat = "#,##0"
$c.Value = 32.0
$c = $ws.Range("H17")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -12.5
$c = $ws.Range("I17")
$c.NumberFormat = "#,##0"
$c.Value = 99.0
$c = $ws.Range("J17")
$c.NumberFormat = "#,##0"
$c.Value = 101.0
$c = $ws.Range("K17")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -1.980198019801
$c = $ws.Range("L17")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 13.793103448275
$c = $ws.Range("M17")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 57.142857142857
$c = $ws.Range("N17")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -65.743944636678

# Row 18
$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D18")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("E18")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -100.0
$c = $ws.Range("F18")
$c.NumberFormat = "#,##0"
$c.Value = 13.0
$c = $ws.Range("G18")
$c.NumberFormat = "#,##0"
$c.Value = 8.0
$c = $ws.Range("H18")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 62.5
$c = $ws.Range("J18")
$c.NumberFormat = "#,##0"
$c.Value = 59.0
$c = $ws.Range("K18")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 1.694915254237
$c = $ws.Range("L18")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 46.341463414634
$c = $ws.Range("M18")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -11.764705882352
$c = $ws.Range("N18")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -91.36690647482

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "#,##0"
$c.Value = 23.0
$c = $ws.Range("E19")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -56.521739130434
$c = $ws.Range("F19")
$c.NumberFormat = "#,##0"
$c.Value = 50.0
$c = $ws.Range("G19")
$c.NumberFormat = "#,##0"
$c.Value = 61.0
$c = $ws.Range("H19")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -18.032786885245
$c = $ws.Range("I19")
$c.NumberFormat = "#,##0"
$c.Value = 181.0
$c = $ws.Range("J19")
$c.NumberFormat = "#,##0"
$c.Value = 222.0
$c = $ws.Range("K19")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -18.468468468468
$c = $ws.Range("L19")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -4.736842105263
$c = $ws.Range("M19")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 60.176991150442
$c = $ws.Range("N19")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -54.292929292929

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "#,##0"
$c.Value = 2.0
$c = $ws.Range("E20")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 250.0
$c = $ws.Range("F20")
$c.NumberFormat = "#,##0"
$c.Value = 27.0
$c = $ws.Range("G20")
$c.NumberFormat = "#,##0"
$c.Value = 16.0
$c = $ws.Range("H20")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 68.75
$c = $ws.Range("I20")
$c.NumberFormat = "#,##0"
$c.Value = 91.0
$c = $ws.Range("J20")
$c.NumberFormat = "#,##0"
$c.Value = 127.0
$c = $ws.Range("K20")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -28.346456692913
$c = $ws.Range("L20")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 42.1875
$c = $ws.Range("M20")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 225.0
$c = $ws.Range("N20")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -87.735849056603

# Row 21
$c = $ws.Range("C21")
$c.NumberFormat = "#,##0"
$c.Value = 22.0
$c = $ws.Range("D21")
$c.NumberFormat = "#,##0"
$c.Value = 37.0
$c = $ws.Range("E21")
$c.NumberFormat = "#,##0.00;""-""#,##0.00"
$c.Value = -40.54054054054
$c = $ws.Range("F21")
$c.NumberFormat = "#,##0"
$c.Value = 133.0
$c = $ws.Range("G21")
$c.NumberFormat = "#,##0"
$c.Value = 131.0
$c = $ws.Range("H21")
$c.NumberFormat = "#,##0.00;""-""#,##0.00"
$c.Value = 1.526717557251
$c = $ws.Range("I21")
$c.NumberFormat = "#,##0"
$c.Value = 515.0
$c = $ws.Range("J21")
$c.NumberFormat = "#,##0"
$c.Value = 595.0
$c = $ws.Range("K21")
$c.NumberFormat = "#,##0.00;""-""#,##0.00"
$c.Value = -13.44537815126
$c = $ws.Range("L21")
$c.NumberFormat = "#,##0.00;""-""#,##0.00"
$c.Value = 14.444444444444
$c = $ws.Range("M21")
$c.NumberFormat = "#,##0.00;""-""#,##0.00"
$c.Value = 33.419689119171
$c = $ws.Range("N21")
$c.NumberFormat = "#,##0.00;""-""#,##0.00"
$c.Value = -80.975249353527

# Row 22
$c = $ws.Range("C22")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("D22")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("E22")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 0.0
$c = $ws.Range("F22")
$c.NumberFormat = "#,##0"
$c.Value = 7.0
$c = $ws.Range("H22")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 133.333333333333
$c = $ws.Range("I22")
$c.NumberFormat = "#,##0"
$c.Value = 18.0
$c = $ws.Range("J22")
$c.NumberFormat = "#,##0"
$c.Value = 14.0
$c = $ws.Range("K22")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 28.571428571428
$c = $ws.Range("L22")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 80.0
$c = $ws.Range("M22")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 100.0

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("E23")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -100.0
$c = $ws.Range("J23")
$c.NumberFormat = "#,##0"
$c.Value = 12.0
$c = $ws.Range("K23")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -8.333333333333
$c = $ws.Range("L23")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 22.222222222222

# Row 24
$c = $ws.Range("C24")
$c.NumberFormat = "#,##0"
$c.Value = 17.0
$c = $ws.Range("D24")
$c.NumberFormat = "#,##0"
$c.Value = 16.0
$c = $ws.Range("E24")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 6.25
$c = $ws.Range("F24")
$c.NumberFormat = "#,##0"
$c.Value = 83.0
$c = $ws.Range("G24")
$c.NumberFormat = "#,##0"
$c.Value = 72.0
$c = $ws.Range("H24")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 15.277777777777
$c = $ws.Range("I24")
$c.NumberFormat = "#,##0"
$c.Value = 421.0
$c = $ws.Range("J24")
$c.NumberFormat = "#,##0"
$c.Value = 455.0
$c = $ws.Range("K24")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -7.472527472527
$c = $ws.Range("L24")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 63.813229571984
$c = $ws.Range("M24")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 109.452736318408

# Row 25
$c = $ws.Range("C25")
$c.NumberFormat = "#,##0"
$c.Value = 10.0
$c = $ws.Range("D25")
$c.NumberFormat = "#,##0"
$c.Value = 9.0
$c = $ws.Range("E25")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 11.111111111111
$c = $ws.Range("F25")
$c.NumberFormat = "#,##0"
$c.Value = 31.0
$c = $ws.Range("H25")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -11.428571428571
$c = $ws.Range("I25")
$c.NumberFormat = "#,##0"
$c.Value = 165.0
$c = $ws.Range("J25")
$c.NumberFormat = "#,##0"
$c.Value = 156.0
$c = $ws.Range("K25")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 5.76923076923
$c = $ws.Range("L25")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 30.952380952381
$c = $ws.Range("M25")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -10.326086956521

# Row 26
$c = $ws.Range("C26")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("F26")
$c.NumberFormat = "#,##0"
$c.Value = 2.0
$c = $ws.Range("H26")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 0.0
$c = $ws.Range("I26")
$c.NumberFormat = "#,##0"
$c.Value = 8.0
$c = $ws.Range("K26")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -46.666666666666
$c = $ws.Range("L26")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -46.666666666666

# Row 27
$c = $ws.Range("C27")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c = $ws.Range("G27")
$c.NumberFormat = "#,##0"
$c.Value = 2.0
$c = $ws.Range("H27")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 150.0
$c = $ws.Range("I27")
$c.NumberFormat = "#,##0"
$c.Value = 21.0
$c = $ws.Range("K27")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 10.526315789473
$c = $ws.Range("L27")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = 5.0

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c = $ws.Range("F28")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("G28")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("M28")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -80.0
$c = $ws.Range("N28")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -98.666666666666

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "***.*"
$c = $ws.Range("F29")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("G29")
$c.NumberFormat = "#,##0"
$c.Value = 1.0
$c = $ws.Range("M29")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -80.0
$c = $ws.Range("N29")
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -98.571428571428

